# Updates cryptos list cell values to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.134.39"
$ws.Range("E2").Value = "  +2.72%  "

$ws.Range("D3").Value = "3.996.01"
$ws.Range("E3").Value = "  +1.36%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "'530.85"
$ws.Range("E5").Value = "  +6.59%  "

$ws.Range("D6").Value = "'147.50"
$ws.Range("E6").Value = "  -0.47%  "

$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  -0.37%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").Value = "'0.742"
$ws.Range("E9").Value = "  +0.93%  "

$ws.Range("D10").Value = "'0.177"
$ws.Range("E10").Value = "  +0.42%  "

$ws.Range("D11").Value = "'0.0000350"
$ws.Range("E11").Value = "  -0.07%  "

$ws.Range("D12").Value = "'42.89"
$ws.Range("E12").Value = "  -1.16%  "

$ws.Range("D13").Value = "'10.53"
$ws.Range("E13").Value = "  +0.38%  "

$ws.Range("D14").Value = "4.589.53"
$ws.Range("E14").Value = "  +0.33%  "

$ws.Range("D15").Value = "3.972.19"
$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("D16").Value = "'21.37"
$ws.Range("E16").Value = "  +6.81%  "

$ws.Range("D17").Value = "'14.27"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("E18").Value = "  +2.69%  "

$ws.Range("D19").Value = "'0.134"
$ws.Range("E19").Value = "  -1.80%  "

$ws.Range("D20").Value = "70.930.14"
$ws.Range("E20").Value = "  +2.32%  "

$ws.Range("D21").Value = "'440.70"
$ws.Range("E21").Value = "  +0.77%  "

$ws.Range("D22").Value = "'3.55"
$ws.Range("E22").Value = "  +2.70%  "

$ws.Range("D23").Value = "'90.35"
$ws.Range("E23").Value = "  +1.60%  "

$ws.Range("D24").Value = "'14.24"
$ws.Range("E24").Value = "  -2.81%  "

$ws.Range("D25").Value = "'4.07"
$ws.Range("E25").Value = "  +5.86%  "

$ws.Range("D26").Value = "'11.85"
$ws.Range("E26").Value = "  -1.67%  "

$ws.Range("D27").Value = "'10.75"
$ws.Range("E27").Value = "  -3.78%  "

$ws.Range("D28").Value = "'37.13"
$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("D29").Value = "'696.92"
$ws.Range("E29").Value = "  -0.46%  "

$ws.Range("D30").Value = "'13.45"
$ws.Range("E30").Value = "  +0.17%  "

$ws.Range("E31").Value = "  -1.45%  "

$ws.Range("D32").Value = "'2.88"
$ws.Range("E32").Value = "  +0.82%  "

$ws.Range("D33").Value = "'6.79"
$ws.Range("E33").Value = "  +12.02%  "

$ws.Range("D34").Value = "'67.39"
$ws.Range("E34").Value = "  +7.61%  "

$ws.Range("D35").Value = "0.0₃0938"
$ws.Range("E35").Value = "  +4.44%  "

$ws.Range("D36").Value = "'0.440"
$ws.Range("E36").Value = "  -3.36%  "

$ws.Range("D37").Value = "'40.32"
$ws.Range("E37").Value = "  -2.00%  "

$ws.Range("D38").Value = "'0.149"
$ws.Range("E38").Value = "  -0.54%  "

$ws.Range("D39").Value = "'3.42"
$ws.Range("E39").Value = "  +10.65%  "

$ws.Range("E40").Value = "  +0.46%  "

$ws.Range("D42").Value = "'0.0486"
$ws.Range("E42").Value = "  -0.40%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'3.17"
$ws.Range("E43").Value = "  +4.73%  "

$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "'2.87"
$ws.Range("E44").Value = "  -1.09%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'3.28"
$ws.Range("E45").Value = "  +9.73%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.53"
$ws.Range("E46").Value = "  +4.34%  "

$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").Value = "'0.000293"
$ws.Range("E47").Value = "  +22.68%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.143"
$ws.Range("E48").Value = "  -0.27%  "

$ws.Range("D49").Value = "'9.24"
$ws.Range("E49").Value = "  +6.04%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0353"
$ws.Range("E50").Value = "  +0.27%  "

$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "'3.38"
$ws.Range("E51").Value = "  -0.15%  "
